$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at 15 for the new "Regresión lineal" task (shifts old rows 15-21 -> 16-22, 24 -> 25) ---
$ws.Rows("15:15").Insert()

# Excel's row-insert carries the formatting of the row above into the leftover
# numeric columns (H:L and AC:AG) that belonged to the "Algoritmo Q-learning" row.
# The new row only has two populated cells (A and AI), so clear those stray cells.
$ws.Range("H15:L15").Clear()
$ws.Range("AC15:AG15").Clear()

# --- New task row content ---
$ws.Range("A15").Value = "Regresión lineal"
$ws.Range("G13").Copy($ws.Range("AI15"))
$ws.Range("AI15").Value = "6 h."

# --- Extend the date header row with one more day, and bump the existing last date ---
$ws.Range("AH2").Copy($ws.Range("AI2"))
$ws.Range("AH2").Value = 44035
$ws.Range("AI2").Value = 44036

# --- Renamed / reworked tasks ---
$ws.Range("A11").Value = "Regresión lineal"
$ws.Range("A12").Value = "Implementación aprendizaje automático"
$ws.Range("A13").Value = "Estructuras básicas Q-learning"

# --- Updated totals (row 24 shifted to row 25 by the insert above) ---
$ws.Range("A25").Value = "4+3+1.5+4+4+3+2+3.5+4.5+4.5+5+1+5+3+4+5+4+4+3+4+3.5+8+3.5+5+5+5+4+3+3.5+4+2.5+5.5+6"
$ws.Range("A1").Value = "Total horas: 130.5"

# --- Update the saved selection to match the edited document ---
$ws.Range("J33").Select()
